# Update countries & provincias Spain
# - Refresh the "Datos actualizados" timestamp.
# - Austria (row 20): refreshed case counts.
# - Rumania's numbers improved enough to overtake Noruega and Australia in
#   the ranking, so rows 31-33 now show Rumania / Noruega / Australia (in
#   that order) instead of Noruega / Australia / Rumania, each carrying the
#   stats for its (possibly shifted) row.
# - Albania (row 97): refreshed case counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp header
$ws.Range("A1").Value = "Datos actualizados a 13 de Abril de 2020 a las 12:22"

# Austria
$ws.Range("B20").Value = 13974
$ws.Range("C20").Value = 29
$ws.Range("E20").Value = 6263

# Row 31 now shows Rumania with its updated totals
$ws.Range("A31").Value = "Rumania"
$ws.Range("B31").Value = 6633
$ws.Range("C31").Value = 333
$ws.Range("D31").Value = 914
$ws.Range("E31").Value = 5401
$ws.Range("F31").Value = 231
$ws.Range("G31").Value = 2
$ws.Range("H31").Value = 318

# Row 32 now shows Noruega (its stats are unchanged, just shifted here)
$ws.Range("A32").Value = "Noruega"
$ws.Range("B32").Value = 6525
$ws.Range("C32").Value = 0
$ws.Range("D32").Value = 32
$ws.Range("E32").Value = 6365
$ws.Range("F32").Value = 59
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 128

# Row 33 now shows Australia (its stats are unchanged, just shifted here)
$ws.Range("A33").Value = "Australia"
$ws.Range("B33").Value = 6359
$ws.Range("C33").Value = 46
$ws.Range("D33").Value = 3494
$ws.Range("E33").Value = 2804
$ws.Range("F33").Value = 79
$ws.Range("G33").Value = 2
$ws.Range("H33").Value = 61

# Albania
$ws.Range("B97").Value = 467
$ws.Range("C97").Value = 21
$ws.Range("D97").Value = 232
$ws.Range("E97").Value = 212
